$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for 2025-11-15 (serial 45976), two station rows appended
# after the existing data (rows 152 and 153).

$row152 = @{
    A = 45976
    B = "四方坪站充电量(kw)"
    C = 603.46799999999996
    D = 1137.6000000000004
    E = 472.62999999999994
    F = 318.36099999999999
    G = 433.39700000000005
    H = 627.99300000000005
    I = 682.64800000000002
    J = 105.188
    K = 81.757999999999996
    L = 152.809
    M = 139.11500000000001
    N = 312.71100000000001
    O = 661.42500000000018
    P = 1323.2009999999998
    Q = 547.1400000000001
    R = 639.99199999999996
    S = 454.41300000000001
    T = 108.733
    U = 144.44
    V = 41.37
    W = 107.88000000000001
    X = 76.16
    Y = 133.64000000000001
    Z = 36.69
}

$row153 = @{
    A = 45976
    B = "高岭站充电量(kw)"
    C = 355.154
    D = 666.49900000000002
    E = 59.892000000000003
    F = 137.85400000000001
    G = 68.26400000000001
    H = 279.78099999999995
    I = 85.421000000000006
    J = 85.734000000000009
    K = 243.41499999999999
    L = 193.952
    M = 82.155000000000001
    N = 187.35899999999998
    O = 273.80100000000004
    P = 378.71399999999994
    Q = 199.51
    R = 220.94400000000002
    S = 104.372
    T = 256.47399999999999
    U = 72.444999999999993
    V = 0
    W = 0
    X = 110.96199999999999
    Y = 51.497999999999998
    Z = 0
}

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

foreach ($col in $cols) {
    $ws.Range($col + "152").Value = $row152[$col]
    $ws.Range($col + "153").Value = $row153[$col]
}

$ws.Range("F158").Select()
